# Update flyer and committee
# Applies four small text edits inside the "TextBox 14" shape that lives
# nested inside the "Group 23" group shape on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the group ("Group 23") and the text box inside it ("TextBox 14")
# by name so the script is resilient to shape ordering.
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Group 23") {
        $grp = $sh
        break
    }
}

$tb = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($item.Name -eq "TextBox 14") {
        $tb = $item
        break
    }
}

$tr = $tb.TextFrame.TextRange

# 1) "@ Toronto " -> "Toronto 2024 "
$idx = $tr.Text.IndexOf("@ Toronto ")
$tr.Characters($idx + 1, 10).Text = "Toronto 2024 "

# 2) ", with " -> ", "
$idx = $tr.Text.IndexOf(", with ")
$tr.Characters($idx + 1, 7).Text = ", "

# 3) "or just with " -> "or just your "
$idx = $tr.Text.IndexOf("or just with ")
$tr.Characters($idx + 1, 13).Text = "or just your "

# 4) "Registration will open soon!" -> "Registration will open soon."
$idx = $tr.Text.IndexOf("Registration will open soon!")
$tr.Characters($idx + 1, 29).Text = "Registration will open soon."
